# Documentation File.docx - final project write-up additions.
#
# The last paragraph ("The server has to be running first or the program
# wont run.") gets re-run into three runs so the misspelling "wont" is
# isolated and flagged with proofErr markers (spellcheck + grammar squiggle
# simulation), and several new paragraphs are appended after it: a note
# about the server's output file, spacer paragraphs, a bold/underlined
# "GitHub Link:" heading, and the repo URL (carrying forward the _GoBack
# bookmark), followed by two trailing blank paragraphs.

$d = $word.ActiveDocument

# Locate the sentence we need to rewrite/extend and grab it (plus its
# paragraph mark, so the replacement can also carry the new paragraphs
# that follow) as a single Range.
$target = $d.Content
$found = $target.Find.Execute(
    "The server has to be running first or the program wont run.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'server has to be running' sentence"
}
$target.MoveEnd(1, 1) | Out-Null   # extend over the paragraph mark

$rPr      = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$rPrBold  = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="36"/><w:szCs w:val="36"/><w:u w:val="single"/></w:rPr>'

# Paragraph 1: the original paragraph, kept (same paraId/rsid), but split
# into three runs with proofErr wrapping around "wont".
$p1 = '<w:p w14:paraId="6EF32715" w14:textId="7D6B23B7" w:rsidR="00D23A4D" w:rsidRPr="00D23A4D" w:rsidRDefault="00D23A4D" w:rsidP="000D07F1">' +
      '<w:pPr>' + $rPr + '</w:pPr>' +
      '<w:r>' + $rPr + '<w:t xml:space="preserve">The server has to be running first or the program </w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
      '<w:r>' + $rPr + '<w:t>wont</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' +
      '<w:r>' + $rPr + '<w:t xml:space="preserve"> run.</w:t></w:r>' +
      '</w:p>'

# Paragraph 2: new sentence about the server's output file ("thee units").
$p2 = '<w:p><w:pPr>' + $rPr + '</w:pPr>' +
      '<w:r>' + $rPr + '<w:t xml:space="preserve">After the running of the program there should be a file in the folder of the server program that holds the output of </w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r>' + $rPr + '<w:t>thee</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r>' + $rPr + '<w:t xml:space="preserve"> units.</w:t></w:r>' +
      '</w:p>'

# Paragraphs 3-4: blank spacer paragraphs.
$p3 = '<w:p><w:pPr>' + $rPr + '</w:pPr></w:p>'
$p4 = '<w:p><w:pPr>' + $rPr + '</w:pPr></w:p>'

# Paragraph 5: bold/underlined "GitHub Link:" heading.
$p5 = '<w:p><w:pPr>' + $rPrBold + '</w:pPr>' +
      '<w:r>' + $rPrBold + '<w:t>GitHub Link</w:t></w:r>' +
      '<w:r>' + $rPrBold + '<w:t>:</w:t></w:r>' +
      '</w:p>'

# Paragraph 6: the repository URL; the _GoBack bookmark moves down to here.
$p6 = '<w:p><w:pPr>' + $rPr + '</w:pPr>' +
      '<w:r>' + $rPr + '<w:t>https://github.com/dwightthomas/Data-Comm-and-Net-Prog</w:t></w:r>' +
      '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
      '</w:p>'

# Paragraphs 7-8: trailing blank paragraphs.
$p7 = '<w:p><w:pPr>' + $rPr + '</w:pPr></w:p>'
$p8 = '<w:p><w:pPr>' + $rPr + '</w:pPr></w:p>'

$body = $p1 + $p2 + $p3 + $p4 + $p5 + $p6 + $p7 + $p8

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $body + '</w:body></w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

# InsertXML replaces the addressed range's contents in one shot, which is
# exactly what we want here: the old single-run sentence + its paragraph
# mark become the 8 paragraphs built above (the _GoBack bookmark is
# reinstated on the new URL paragraph as part of that payload).
$target.InsertXML($xml) | Out-Null
